$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4: plain numeric data
$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 10

$ws.Range("A3").Value = 1000
$ws.Range("B3").Value = 2000
$ws.Range("C3").Value = 10

$ws.Range("A4").Value = 1000
$ws.Range("B4").Value = 3000
$ws.Range("C4").Value = 20

# Row 5: values entered as text (leading apostrophe forces text, matching
# how the source workbook stores these as strings rather than numbers)
$ws.Range("A5").Value = "'5000.0"
$ws.Range("B5").Value = "'1000.0"
$ws.Range("C5").Value = "'10"
